$d = $word.ActiveDocument

# The Pearson logo inline pictures (footer1.xml / footer2.xml) were renamed
# from image2.png -> image1.png, and the BTec logo inline picture
# (header1.xml) was renamed from image1.jpg -> image2.jpg. Word's
# InlineShape object does not expose a settable Name property (that only
# exists on the floating Shape object), so we round-trip the document's
# WordOpenXML, rewrite the <wp:docPr>/<pic:cNvPr> name="..." attributes in
# place, and hand the edited markup back to Word.

$xml = $d.WordOpenXML

$xml = $xml -replace 'name="image2\.png"', 'name="image1.png"'
$xml = $xml -replace 'name="image1\.jpg"', 'name="image2.jpg"'

$d.WordOpenXML = $xml

Write-Output "renamed inline picture parts"
